$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Updated "Ore lavorate" (B) and "Produzione" (C) values per row.
# Column E holds a shared formula (C-(C*D)/100) that recalculates automatically.
$updates = @(
    @{ Row = 2;  B = 96;  C = 34.92 },
    @{ Row = 3;  B = 88;  C = 39.72 },
    @{ Row = 5;  B = 96;  C = 34.9 },
    @{ Row = 6;  B = 96;  C = 35.28 },
    @{ Row = 7;  C = 33.04 },
    @{ Row = 8;  B = 40;  C = 16.07 },
    @{ Row = 9;  C = 30.3 },
    @{ Row = 10; B = 68;  C = 28.47 },
    @{ Row = 11; B = 96;  C = 36.5 },
    @{ Row = 12; B = 96;  C = 33.9 },
    @{ Row = 13; B = 96;  C = 37.16 },
    @{ Row = 14; B = 96;  C = 35.33 },
    @{ Row = 16; B = 104; C = 37.16 },
    @{ Row = 18; B = 104; C = 26.51 },
    @{ Row = 19; C = 36.36 },
    @{ Row = 20; B = 72;  C = 37.37 },
    @{ Row = 21; B = 96;  C = 32.41 },
    @{ Row = 22; C = 34.9 },
    @{ Row = 23; C = 33.79 },
    @{ Row = 24; B = 96;  C = 36.53 },
    @{ Row = 25; B = 88;  C = 35.85 }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($u.Row, 2).Value = $u.B
    }
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}

# Update the view: zoom to 110% and move the active selection to E22.
$ws.Activate()
$excel.ActiveWindow.Zoom = 110
$ws.Range("E22").Select()
